# ISIS1225 - Tablas de Datos Lab 7
#
# The "Carga de Catalogo PROBING" table (Table1, A2:C6) loses its third
# data row (Factor de Carga 0.6 / 110 / 65) -- the row below it
# (0.8 / 115 / 75) takes its place, shrinking the table to A2:C5.
# The first remaining data row's "Factor de Carga" value also changes
# from 0.2 to 0.3.
#
# The "Carga de Catalogo CHAINING" table (Table13, originally A10:C14)
# shifts up one row (because of the deletion above) and loses its last
# data row (Factor de Carga 8 / 160 / 120), shrinking it to A9:C12.
#
# Deleting whole rows lets Excel re-flow the merged title cells, the
# table definitions (ref/autoFilter) and the sheet dimension on its own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab7")

# Drop the PROBING row "0.6 / 110 / 65" (row 5) -- row 6 (0.8/115/75)
# slides up into row 5, and everything below shifts up by one row.
$ws.Rows("5:5").Delete()

# First PROBING data row's load factor changes from 0.2 to 0.3.
$ws.Range("A3").Value2 = 0.3

# Drop the last CHAINING row "8 / 160 / 120". After the deletion above
# it now lives at row 13 (was row 14).
$ws.Rows("13:13").Delete()

# Reflect the author's final selection on the data sheet.
$ws.Range("B25").Select()
